# Applies the "Updated cryptos list" data refresh to the crypto-tracker sheet.
# Column D ("Price") values are numeric-looking text (e.g. "417.70", "1.00") that
# must stay as literal text (exact digits/trailing zeros), so we force a Text
# number format before assigning, then restore the default "Normal" style so we
# do not leave a stray number format applied to the cell.
# Column E ("Volume(1h)") values already contain padding spaces, so Excel keeps
# them as text automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.197.48'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +5.42%  '
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.507.25'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.71%  '
# Row 4
$ws.Range("E4").Value = '  +0.02%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '417.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.00%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.59'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.88%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.660'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.70%  '
# Row 8
$ws.Range("E8").Value = '  +0.03%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.783'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.62%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.163'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +14.82%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '43.59'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.21%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000264'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +18.45%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.96'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +8.54%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.059.15'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.69%  '
# Row 15
$ws.Range("E15").Value = '  +0.57%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.48'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.17%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.509.63'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.24%  '
# Row 18
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.75'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.09%  '
# Row 19
$ws.Range("B19").Value = 'Polygon'
$ws.Range("C19").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.10'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.98%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '65.040.47'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.10%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '459.06'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.06%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '90.54'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.77%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.21'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.95%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.34'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.00%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.38'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.52%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.93'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.06%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '34.15'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.74%  '
# Row 28
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.53'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.99%  '
# Row 29
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.72'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.94%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.50'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.40%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.116'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.86%  '
# Row 32
$ws.Range("E32").Value = '  -1.62%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '39.82'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.75%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.19%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '57.38'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.99%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0506'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.04%  '
# Row 37
$ws.Range("B37").Value = 'Stellar'
$ws.Range("C37").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.153'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +13.53%  '
# Row 38
$ws.Range("B38").Value = 'PEPE'
$ws.Range("C38").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0712'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +34.49%  '
# Row 39
$ws.Range("E39").Value = '  +2.45%  '
# Row 40
$ws.Range("E40").Value = '  -0.11%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.78'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.98%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.55'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.72%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '145.50'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.00%  '
# Row 44
$ws.Range("E44").Value = '  -0.93%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.313'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.96%  '
# Row 46
$ws.Range("E46").Value = '  -3.80%  '
# Row 47
$ws.Range("E47").Value = '  -0.22%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.85'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.22%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.144'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.07%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '21.69'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.77%  '
# Row 51
$ws.Range("B51").Value = 'ApeXProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.47'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.40%  '
